$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text storage for the Price column (values are stored as text, not numbers)
$cells = @{
    "D2" = "241.76"
    "D3" = "21.75"
    "D4" = "5.371"
    "D5" = "0.05688"
    "D6" = "3.415"
    "D7" = "6.295"
    "D8" = "0.8060"
    "D9" = "0.8491"
    "D11" = "0.07258"
    "D12" = "0.03052"
    "D13" = "0.03141"
    "D14" = "0.09356"
    "D15" = "3.922"
    "D16" = "0.001578"
    "D17" = "0.04814"
    "D18" = "0.0005825"
    "D19" = "0.006393"
    "D20" = "0.004058"
    "D21" = "0.0009927"
    "D22" = "0.0001501"
    "D23" = "3.736"
    "D24" = "2.151"
    "D25" = "0.3260"
    "D26" = "0.1311"
    "D27" = "0.0004002"
    "D40" = "0.03809"
    "D41" = "0.006702"
    "D43" = "0.002622"
    "D44" = "0.006878"
    "D45" = "0.00005624"
    "D46" = "0.00000000751"
    "D47" = "0.5805"
    "D48" = "0.1423"
    "D49" = "0.00002102"
    "D50" = "0.01011"
}

foreach ($addr in $cells.Keys) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $cells[$addr]
}

